# Update "paises.xlsx" (Pais sheet): refresh COVID country stats + reorder
# two pairs of countries in the underlying string table (Estonia/Ruanda and
# Islas Malvinas/Montserrat swap positions) + bump the "datos actualizados"
# timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp banner (row 1) ---------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 12 de Agosto de 2020 a las 10:46"

# --- Country name reorder (Estonia <-> Ruanda) -----------------------------
$ws.Range("A130").Value = "Estonia"
$ws.Range("A131").Value = "Ruanda"

# --- Country name reorder (Islas Malvinas <-> Montserrat) ------------------
$ws.Range("A213").Value = "Islas Malvinas"
$ws.Range("A214").Value = "Montserrat"

# --- Updated case counts -----------------------------------------------
# Row 7 - Rusia
$ws.Range("B7").Value = 902701
$ws.Range("C7").Value = 5102
$ws.Range("D7").Value = 710298
$ws.Range("E7").Value = 177143
$ws.Range("G7").Value = 129
$ws.Range("H7").Value = 15260

# Row 25 - Filipinas
$ws.Range("B25").Value = 143749
$ws.Range("C25").Value = 4444
$ws.Range("D25").Value = 68997
$ws.Range("E25").Value = 72348
$ws.Range("G25").Value = 93
$ws.Range("H25").Value = 2404

# Row 26 - Indonesia
$ws.Range("B26").Value = 130718
$ws.Range("C26").Value = 1942
$ws.Range("D26").Value = 85798
$ws.Range("E26").Value = 39017
$ws.Range("G26").Value = 79
$ws.Range("H26").Value = 5903

# Row 47 - Singapur
$ws.Range("B47").Value = 55395
$ws.Range("C47").Value = 42
$ws.Range("E47").Value = 5240

# Row 48 - Polonia
$ws.Range("B48").Value = 53676
$ws.Range("C48").Value = 715
$ws.Range("D48").Value = 37611
$ws.Range("E48").Value = 14235
$ws.Range("G48").Value = 9
$ws.Range("H48").Value = 1830

# Row 71 - Austria
$ws.Range("B71").Value = 22439
$ws.Range("C71").Value = 194
$ws.Range("D71").Value = 20268
$ws.Range("E71").Value = 1447
$ws.Range("G71").Value = 1
$ws.Range("H71").Value = 724

# Row 111 - Hong Kong
$ws.Range("E111").Value = 1069
$ws.Range("G111").Value = 3
$ws.Range("H111").Value = 61

# Row 123 - Eslovaquia
$ws.Range("B123").Value = 2690
$ws.Range("C123").Value = 75
$ws.Range("D123").Value = 1884
$ws.Range("E123").Value = 775

# Row 128 - Lituania
$ws.Range("B128").Value = 2309
$ws.Range("C128").Value = 26
$ws.Range("D128").Value = 1683
$ws.Range("E128").Value = 545

# Row 130 - Estonia (data now attached to this row after the name swap)
$ws.Range("B130").Value = 2174
$ws.Range("C130").Value = 7
$ws.Range("D130").Value = 1975
$ws.Range("E130").Value = 136
$ws.Range("H130").Value = 63

# Row 131 - Ruanda (data now attached to this row after the name swap)
$ws.Range("B131").Value = 2171
$ws.Range("D131").Value = 1478
$ws.Range("E131").Value = 686
$ws.Range("H131").Value = 7

# Row 159 - Bahamas
$ws.Range("D159").Value = 400
$ws.Range("E159").Value = 449

# Row 174 - Islas Feroe
$ws.Range("B174").Value = 324
$ws.Range("C174").Value = 6
$ws.Range("E174").Value = 99

# Row 213 - Islas Malvinas (data now attached to this row after the name swap)
$ws.Range("D213").Value = 13
$ws.Range("H213").Value = 0

# Row 214 - Montserrat (data now attached to this row after the name swap)
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1
